$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 463 (the alpaca post), shifting subsequent rows up.
$ws.Rows.Item(463).Delete()
